# New Data is added: populate the "StudentForm" sheet with the student
# registration data, rename Sheet2 -> StudentForm and make it the active tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Rename Sheet2 -> StudentForm -----------------------------------------
$ws2.Name = "StudentForm"

# --- Row 1: plain column-index markers (0..20) across A1:U1 ---------------
for ($c = 1; $c -le 21; $c++) {
    $ws2.Cells.Item(1, $c).Value = $c - 1
}

# --- Row 2: single marker cell ---------------------------------------------
$ws2.Range("A2").Value = "StudentForm"

# --- Row 3: header labels ----------------------------------------------------
$ws2.Range("A3").Value = "RunMode"
$ws2.Range("B3").Value = "Lastname"
$ws2.Range("C3").Value = "Firstname"
$ws2.Range("D3").Value = "Middlename"
$ws2.Range("E3").Value = "SSN"
$ws2.Range("F3").Value = "Grade"
$ws2.Range("G3").Value = "Studentid"
$ws2.Range("H3").Value = "DOB"
$ws2.Range("I3").Value = "Address"
$ws2.Range("J3").Value = "Zip"
$ws2.Range("K3").Value = "Home"
$ws2.Range("L3").Value = "Cell"
$ws2.Range("M3").Value = "Emailaddress"
$ws2.Range("N3").Value = "Parentemployer"
$ws2.Range("O3").Value = "Telephone"
$ws2.Range("P3").Value = "Name"
$ws2.Range("Q3").Value = "Relationtochild"
$ws2.Range("R3").Value = "Pickuppersonname"
$ws2.Range("S3").Value = "PickuppersonTelephone"
$ws2.Range("T3").Value = "ParentorGuardian"
$ws2.Range("U3").Value = "SignatureLegalParentorGuardian"

# --- Row 4: data values -----------------------------------------------------
$ws2.Range("A4").Value = "StudentForm"

# B4 "Jack" - plain (non-hyperlink) black Calibri, explicitly set so the
# theme/minor-scheme link is dropped (matches an "explicit" font entry).
$ws2.Range("B4").Value = "Jack"
$ws2.Range("B4").Font.Name = "Calibri"

# C4 "Joshua " - right aligned (reuses the existing right-aligned style).
$ws2.Range("C4").Value = "Joshua "
$ws2.Range("C4").HorizontalAlignment = -4152   # xlRight

$ws2.Range("D4").Value = "Jhon"
$ws2.Range("E4").Value = "979-87-8787"
$ws2.Range("F4").Value = "K4"
$ws2.Range("G4").Value = 912

# H4 DOB - a date value, formatted with the built-in short-date format.
$ws2.Range("H4").Value = 38880
$ws2.Range("H4").NumberFormat = "mm-dd-yy"

# I4 Address - smaller black font.
$ws2.Range("I4").Value = "3217 Versante Drive"
$ws2.Range("I4").Font.Size = 9
$ws2.Range("I4").Font.Color = 0

$ws2.Range("J4").Value = 38119
$ws2.Range("K4").Value = "965-679-8989"
$ws2.Range("L4").Value = "966-659-7666"

# M4 Emailaddress - hyperlinked to the same mailto address used on Login.
$ws2.Range("M4").Value = "mani6747@gmail.com"
$ws2.Hyperlinks.Add($ws2.Range("M4"), "mailto:mani6747@gmail.com")

$ws2.Range("N4").Value = "IT Industry"
$ws2.Range("O4").Value = "965-679-8989"
$ws2.Range("P4").Value = "Daren"
$ws2.Range("Q4").Value = "Uncel"
$ws2.Range("R4").Value = "Josep"
$ws2.Range("S4").Value = "965-679-8989"
$ws2.Range("T4").Value = "Parent"
$ws2.Range("U4").Value = "Manikanta"

# --- Column widths: autofit to content --------------------------------------
$ws2.Range("A1:U4").Columns.AutoFit()

# --- Page setup: portrait -----------------------------------------------------
$ws2.PageSetup.Orientation = 1   # xlPortrait

# --- Selection / active sheet bookkeeping ------------------------------------
$ws1.Activate()
$ws1.Range("A1:C4").Select()

$ws2.Activate()
$ws2.Range("L3").Select()

# --- Workbook OLE size (matches the now-larger used range on StudentForm) ---
$wb.Names.Add("_xlnm.OLE_Size", "=StudentForm!`$A`$1:`$L`$12") | Out-Null
